$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "25.966.05"
Set-TextValue $ws "E2" "  -0.34%  "

Set-TextValue $ws "D3" "1.635.07"
Set-TextValue $ws "E3" "  -0.72%  "

Set-TextValue $ws "E4" "  -0.20%  "

Set-TextValue $ws "D5" "214.19"
Set-TextValue $ws "E5" "  -1.19%  "

Set-TextValue $ws "E6" "  -0.72%  "

Set-TextValue $ws "E7" "  -0.09%  "

Set-TextValue $ws "D8" "0.252"
Set-TextValue $ws "E8" "  -1.92%  "

Set-TextValue $ws "E9" "  -2.93%  "

Set-TextValue $ws "D10" "18.51"
Set-TextValue $ws "E10" "  -6.07%  "

Set-TextValue $ws "D11" "0.0791"
Set-TextValue $ws "E11" "  -0.57%  "

Set-TextValue $ws "D12" "1.862.01"
Set-TextValue $ws "E12" "  -0.66%  "

Set-TextValue $ws "D13" "1.642.37"
Set-TextValue $ws "E13" "  -0.30%  "

Set-TextValue $ws "E14" "  -2.96%  "

Set-TextValue $ws "B16" "ShibaInu"
Set-TextValue $ws "C16" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws "D16" "0.0₃0745"
Set-TextValue $ws "E16" "  -2.98%  "

Set-TextValue $ws "B17" "WrappedBTC"
Set-TextValue $ws "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws "D17" "25.981.09"
Set-TextValue $ws "E17" "  -0.18%  "

Set-TextValue $ws "D18" "61.67"
Set-TextValue $ws "E18" "  -2.46%  "

Set-TextValue $ws "E19" "  -0.18%  "

Set-TextValue $ws "D20" "190.54"
Set-TextValue $ws "E20" "  -1.52%  "

Set-TextValue $ws "E21" "  -2.81%  "

Set-TextValue $ws "E22" "  -3.89%  "

Set-TextValue $ws "E23" "  -2.11%  "

Set-TextValue $ws "E24" "  -0.93%  "

Set-TextValue $ws "E25" "  -0.80%  "

Set-TextValue $ws "E26" "  -0.21%  "

Set-TextValue $ws "E27" "  -3.59%  "

Set-TextValue $ws "D28" "6.79"
Set-TextValue $ws "E28" "  -2.12%  "

Set-TextValue $ws "D29" "15.26"
Set-TextValue $ws "E29" "  -1.99%  "

Set-TextValue $ws "D30" "1.23"
Set-TextValue $ws "E30" "  -1.45%  "

Set-TextValue $ws "E31" "  -3.53%  "

Set-TextValue $ws "E32" "  -3.21%  "

Set-TextValue $ws "D33" "3.14"
Set-TextValue $ws "E33" "  -4.66%  "

Set-TextValue $ws "E34" "  -2.07%  "

Set-TextValue $ws "E35" "  -2.71%  "

Set-TextValue $ws "D36" "1.137.51"
Set-TextValue $ws "E36" "  +0.37%  "

Set-TextValue $ws "E37" "  -4.74%  "

Set-TextValue $ws "D38" "2.43"
Set-TextValue $ws "E38" "  -1.47%  "

Set-TextValue $ws "D39" "0.523"
Set-TextValue $ws "E39" "  -3.73%  "

Set-TextValue $ws "E40" "  -1.47%  "

Set-TextValue $ws "D41" "98.50"
Set-TextValue $ws "E41" "  -1.22%  "

Set-TextValue $ws "D42" "0.779"
Set-TextValue $ws "E42" "  -2.34%  "

Set-TextValue $ws "D43" "1.771.96"
Set-TextValue $ws "E43" "  -0.61%  "

Set-TextValue $ws "E44" "  -5.14%  "

Set-TextValue $ws "E45" "  -1.52%  "

Set-TextValue $ws "D46" "55.23"
Set-TextValue $ws "E46" "  -2.75%  "

Set-TextValue $ws "D47" "0.0528"
Set-TextValue $ws "E47" "  -0.37%  "

Set-TextValue $ws "E48" "  +1.04%  "

Set-TextValue $ws "E49" "  -0.62%  "

Set-TextValue $ws "D50" "7.57"
Set-TextValue $ws "E50" "  -1.71%  "

Set-TextValue $ws "E51" "  +0.08%  "
